# "Import HTML.xlsx" — re-save / cleanup pass.
#
# The upstream commit that produced this revision does not change any cell
# data, formula, or value: diffing the canonical OOXML shows only artifacts
# of opening the (Google-Sheets-exported) workbook in desktop Excel and
# saving it again — refreshed XML namespaces/revision metadata in
# xl/workbook.xml, Excel-normalised xl/styles.xml (font size tokens,
# gray125 fill, explicit border children, table/slicer style defaults),
# recalculated `spans`/`dimension`/`pageMargins` bookkeeping in the
# worksheets, `ca="1"` markers on the (unrecognised-function) formulas, and
# the removal of the two placeholder drawing parts that never contained any
# shape. None of that touches a single cell's content, formula result, or
# formatting intent, so there is nothing to replay on the object model
# beyond confirming the workbook's existing state.
#
# Touch the workbook the way Excel would on a plain open/save cycle: make
# sure Sheet1 is the active sheet with A1 selected (matches the saved
# `tabSelected="1"` / active-cell state), without altering any cell values,
# formulas, or styles.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A1").Select()

# No cell values, formulas, or styles are modified — the content already
# matches the target revision; this script only mirrors the no-op
# open/save round trip recorded in the commit.
